# Ascend & Descend Movement Fix
# - Fixed so that Ascend & Descend Movement also works with Slabs and WaterBlocks
#
# The commit mostly consists of cosmetic proofing-run splits (w:proofErr
# gramStart/gramEnd/spellStart/spellEnd markers inserted by Word's grammar/
# spell checker) plus a couple of run merges, and the removal of a now
# redundant block of "To-Do" bullet points that covered Ascend/Descend-only
# behaviour (the fix now generalizes them, so the old itemised notes about
# restart/save-state behaviour etc. were trimmed away).

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$d = $word.ActiveDocument

function Replace-ParagraphXml($index, $innerXml) {
    $p = $d.Paragraphs.Item($index)
    $rng = $p.Range
    $xml = "<w:p $ns>$innerXml</w:p>"
    $rng.InsertXML($xml)
}

# Locate the paragraphs we need to touch by a distinctive substring so the
# script is resilient to minor index drift.
function Find-ParagraphIndex($needle) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        if ($paras.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) "Dark Levels" contain a different/more rare collectable and no coins
#    -> split "more rare" out with gramStart/gramEnd proofErr markers
# ---------------------------------------------------------------------
$i1 = Find-ParagraphIndex("different/more rare")
$inner1 = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='2'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>&#8220;Dark Levels&#8221; contain a different/</w:t></w:r>" +
  "<w:proofErr w:type='gramStart'/>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>more rare</w:t></w:r>" +
  "<w:proofErr w:type='gramEnd'/>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> collectable and no coins</w:t></w:r>"
Replace-ParagraphXml $i1 $inner1

# ---------------------------------------------------------------------
# 2) "Weak Blocks " + "may hide Refill blocks..." -> merge into one run
# ---------------------------------------------------------------------
$i2 = Find-ParagraphIndex("Weak Blocks may hide")
$inner2 = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='3'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Weak Blocks may hide Refill blocks, Teleporters and other blocks under themselves</w:t></w:r>"
Replace-ParagraphXml $i2 $inner2

# ---------------------------------------------------------------------
# 3) "As long as the player moves..." -> split "As long as" with gramStart/gramEnd
# ---------------------------------------------------------------------
$i3 = Find-ParagraphIndex("As long as the player moves")
$inner3 = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='3'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
  "<w:proofErr w:type='gramStart'/>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>As long as</w:t></w:r>" +
  "<w:proofErr w:type='gramEnd'/>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> the player moves in the same direction after stepping on this, it doesn&#8217;t pay any step cost</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> for its movement</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>s</w:t></w:r>"
Replace-ParagraphXml $i3 $inner3

# ---------------------------------------------------------------------
# 4) "Step-on " + "Elevator" -> merge into one run
# ---------------------------------------------------------------------
$i4 = Find-ParagraphIndex("Step-on")
$inner4 = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='2'/><w:numId w:val='1'/></w:numPr><w:rPr><w:color w:val='3A7C22' w:themeColor='accent6' w:themeShade='BF'/><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:color w:val='3A7C22' w:themeColor='accent6' w:themeShade='BF'/><w:lang w:val='en-US'/></w:rPr><w:t>Step-on Elevator</w:t></w:r>"
Replace-ParagraphXml $i4 $inner4

# ---------------------------------------------------------------------
# 5) " (ex. from island to island" -> split "ex" with gramStart/gramEnd
# ---------------------------------------------------------------------
$i5 = Find-ParagraphIndex("ex. from island")
$inner5 = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='3'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Can move over the void</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (</w:t></w:r>" +
  "<w:proofErr w:type='gramStart'/>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>ex</w:t></w:r>" +
  "<w:proofErr w:type='gramEnd'/>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>. from island to island</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> or </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>from </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>one high</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> point to another</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>)</w:t></w:r>"
Replace-ParagraphXml $i5 $inner5

# ---------------------------------------------------------------------
# 6) "...its standingOnBlock and into the Water Block right under (if it
#    can Swim)" -> wrap "standingOnBlock" in spellStart/spellEnd (moving the
#    trailing space out of the "ngOnBlock " run into its own run) and merge
#    the trailing " " + "(if it can Swim)" runs into one.
# ---------------------------------------------------------------------
$i6 = Find-ParagraphIndex("ngOnBlock")
$inner6 = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>Make it so that a player can Descend </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>down </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>through </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>its </w:t></w:r>" +
  "<w:proofErr w:type='spellStart'/>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>stand</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>i</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>ngOnBlock</w:t></w:r>" +
  "<w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>and </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>into </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>the</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> Water Block </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>right under</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (if it can Swim)</w:t></w:r>"
Replace-ParagraphXml $i6 $inner6

# ---------------------------------------------------------------------
# 7) Remove the now-redundant trailing To-Do bullet points: "Restart Map by
#    pressing X" through "Apply saved Data to MainMenu and the Maps when
#    running the scenes" (six paragraphs).
# ---------------------------------------------------------------------
$iStart = Find-ParagraphIndex("Restart Map by pressing X")
$iEnd = Find-ParagraphIndex("Apply saved Data to MainMenu")
$startPos = $d.Paragraphs.Item($iStart).Range.Start
$endPos = $d.Paragraphs.Item($iEnd).Range.End
$delRng = $d.Range($startPos, $endPos)
$delRng.Delete()
